# Auto-generated Excel COM-interop script to apply the numeric updates
# described by the commit diff, across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 34728.062
$ws.Range("J62").Value = 5224.7
$ws.Range("L62").Value = 5224.7
$ws.Range("N62").Value = -6472.7
$ws.Range("H65").Value = 34728.062
$ws.Range("J65").Value = 5224.7
$ws.Range("L65").Value = 26123.5
$ws.Range("N65").Value = -32363.5
$ws.Range("H92").Value = 630.913
$ws.Range("I92").Value = 523.2727
$ws.Range("J92").Value = 2999
$ws.Range("K92").Value = 523.2727
$ws.Range("L92").Value = 2999
$ws.Range("M92").Value = 724.7273
$ws.Range("N92").Value = -5495
$ws.Range("H116").Value = 4659.5625
$ws.Range("I116").Value = 2982.2856
$ws.Range("J116").Value = 5964.1113
$ws.Range("K116").Value = 2982.2856
$ws.Range("L116").Value = 5964.1113
$ws.Range("M116").Value = 459.7143999999998
$ws.Range("N116").Value = -12848.1113
$ws.Range("H127").Value = 5074.8887
$ws.Range("I127").Value = 5074.8887
$ws.Range("K127").Value = 15224.6661
$ws.Range("M127").Value = -10264.6661
$ws.Range("H132").Value = 9850.916999999999
$ws.Range("I132").Value = 8561.286
$ws.Range("K132").Value = 25683.858
$ws.Range("M132").Value = -23153.858
$ws.Range("H137").Value = 3168.0322
$ws.Range("I137").Value = 2792.1765
$ws.Range("K137").Value = 8376.529500000001
$ws.Range("M137").Value = -5826.529500000001
$ws.Range("H138").Value = 167249.27
$ws.Range("I138").Value = 1587.5714
$ws.Range("J138").Value = 254221.66
$ws.Range("K138").Value = 4762.7142
$ws.Range("L138").Value = 762664.98
$ws.Range("M138").Value = 377.2857999999997
$ws.Range("N138").Value = -772944.98

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8993.450999999999
$ws.Range("I32").Value = 8993.450999999999
$ws.Range("K32").Value = 8993.450999999999
$ws.Range("M32").Value = -8706.450999999999
$ws.Range("H34").Value = 34012.5
$ws.Range("I34").Value = 34012.5
$ws.Range("K34").Value = 34012.5
$ws.Range("M34").Value = -33741.5
$ws.Range("H45").Value = 3907.875
$ws.Range("I45").Value = 2635.3635
$ws.Range("J45").Value = 4984.615
$ws.Range("K45").Value = 2635.3635
$ws.Range("L45").Value = 4984.615
$ws.Range("M45").Value = -2258.3635
$ws.Range("N45").Value = -5738.615
$ws.Range("H97").Value = 1183.4
$ws.Range("I97").Value = 966.7273
$ws.Range("K97").Value = 966.7273
$ws.Range("M97").Value = -470.7273
$ws.Range("H122").Value = 2007.0834
$ws.Range("I122").Value = 1915.9032
$ws.Range("K122").Value = 5747.7096
$ws.Range("M122").Value = -3297.7096

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2494
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2494
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2494
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -2774
$ws.Range("H22").Value = 553.44446
$ws.Range("I22").Value = 570.1429000000001
$ws.Range("K22").Value = 570.1429000000001
$ws.Range("M22").Value = -397.1429000000001
$ws.Range("H134").Value = 1576.4849
$ws.Range("J134").Value = 4443
$ws.Range("L134").Value = 13329
$ws.Range("N134").Value = -18399

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2983.8572
$ws.Range("I31").Value = 2284.3333
$ws.Range("J31").Value = 3916.5557
$ws.Range("K31").Value = 2284.3333
$ws.Range("L31").Value = 3916.5557
$ws.Range("M31").Value = -1989.3333
$ws.Range("N31").Value = -4506.5557
$ws.Range("H34").Value = 2983.8572
$ws.Range("I34").Value = 2284.3333
$ws.Range("J34").Value = 3916.5557
$ws.Range("K34").Value = 2284.3333
$ws.Range("L34").Value = 3916.5557
$ws.Range("M34").Value = -2082.3333
$ws.Range("N34").Value = -4320.5557
$ws.Range("H58").Value = 5135.273
$ws.Range("I58").Value = 5134.5557
$ws.Range("J58").Value = 5138.5
$ws.Range("K58").Value = 5134.5557
$ws.Range("L58").Value = 5138.5
$ws.Range("M58").Value = -4931.5557
$ws.Range("N58").Value = -5544.5
$ws.Range("H86").Value = 7450
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 7450
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 9241.071
$ws.Range("I99").Value = 9833.1
$ws.Range("K99").Value = 9833.1
$ws.Range("M99").Value = -8335.1
$ws.Range("H107").Value = 602.13794
$ws.Range("I107").Value = 438.8
$ws.Range("J107").Value = 777.1429000000001
$ws.Range("K107").Value = 438.8
$ws.Range("L107").Value = 777.1429000000001
$ws.Range("M107").Value = 1481.2
$ws.Range("N107").Value = -4617.1429
$ws.Range("H126").Value = 9241.071
$ws.Range("I126").Value = 9833.1
$ws.Range("K126").Value = 29499.3
$ws.Range("M126").Value = -27029.3
$ws.Range("H136").Value = 5135.273
$ws.Range("I136").Value = 5134.5557
$ws.Range("J136").Value = 5138.5
$ws.Range("K136").Value = 15403.6671
$ws.Range("L136").Value = 15415.5
$ws.Range("M136").Value = -12853.6671
$ws.Range("N136").Value = -20515.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 149.88889
$ws.Range("J2").Value = 118.63636
$ws.Range("L2").Value = 711.81816
$ws.Range("N2").Value = -937.81816
$ws.Range("H14").Value = 3325.0833
$ws.Range("I14").Value = 3325.0833
$ws.Range("K14").Value = 9975.249899999999
$ws.Range("M14").Value = -9802.249899999999
$ws.Range("H38").Value = 2544.5833
$ws.Range("I38").Value = 246.75
$ws.Range("J38").Value = 4842.4165
$ws.Range("K38").Value = 740.25
$ws.Range("L38").Value = 14527.2495
$ws.Range("M38").Value = -393.25
$ws.Range("N38").Value = -15221.2495
$ws.Range("H122").Value = 903
$ws.Range("J122").Value = 1049.8334
$ws.Range("L122").Value = 9448.500599999999
$ws.Range("N122").Value = -14348.5006
$ws.Range("H131").Value = 1091
$ws.Range("I131").Value = 824.36365
$ws.Range("J131").Value = 1677.6
$ws.Range("K131").Value = 2473.09095
$ws.Range("L131").Value = 5032.799999999999
$ws.Range("M131").Value = 2566.90905
$ws.Range("N131").Value = -15112.8
$ws.Range("H137").Value = 4539.4707
$ws.Range("I137").Value = 1735.125
$ws.Range("J137").Value = 7032.222
$ws.Range("K137").Value = 5205.375
$ws.Range("L137").Value = 21096.666
$ws.Range("M137").Value = -105.375
$ws.Range("N137").Value = -31296.666
$ws.Range("H139").Value = 3268.7827
$ws.Range("I139").Value = 1922.9
$ws.Range("K139").Value = 5768.700000000001
$ws.Range("M139").Value = -628.7000000000007

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 15009667
$ws.Range("I24").Value = 30006002
$ws.Range("J24").Value = 13332.667
$ws.Range("K24").Value = 30006002
$ws.Range("L24").Value = 13332.667
$ws.Range("M24").Value = -30005829
$ws.Range("N24").Value = -13678.667
$ws.Range("H97").Value = 3045.5454
$ws.Range("I97").Value = 4155.8335
$ws.Range("K97").Value = 4155.8335
$ws.Range("M97").Value = -3659.8335
$ws.Range("H126").Value = 7342.5713
$ws.Range("I126").Value = 7342.5713
$ws.Range("K126").Value = 22027.7139
$ws.Range("M126").Value = -19557.7139
$ws.Range("H132").Value = 8944
$ws.Range("I132").Value = 9220.368
$ws.Range("K132").Value = 27661.104
$ws.Range("M132").Value = -25131.104

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15500
$ws.Range("J5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("N5").Value = -30226
$ws.Range("H16").Value = 680.1
$ws.Range("I16").Value = 265.44446
$ws.Range("J16").Value = 4412
$ws.Range("K16").Value = 265.44446
$ws.Range("L16").Value = 4412
$ws.Range("M16").Value = -95.44445999999999
$ws.Range("N16").Value = -4752
$ws.Range("H40").Value = 3225.9788
$ws.Range("I40").Value = 3225.9788
$ws.Range("K40").Value = 3225.9788
$ws.Range("M40").Value = -3089.9788
$ws.Range("H82").Value = 9930.65
$ws.Range("I82").Value = 9890.933999999999
$ws.Range("J82").Value = 10049.8
$ws.Range("K82").Value = 9890.933999999999
$ws.Range("L82").Value = 10049.8
$ws.Range("M82").Value = -9529.933999999999
$ws.Range("N82").Value = -10771.8
$ws.Range("H85").Value = 9930.65
$ws.Range("I85").Value = 9890.933999999999
$ws.Range("J85").Value = 10049.8
$ws.Range("K85").Value = 9890.933999999999
$ws.Range("L85").Value = 10049.8
$ws.Range("M85").Value = -8642.933999999999
$ws.Range("N85").Value = -12545.8
$ws.Range("H122").Value = 5802.517
$ws.Range("I122").Value = 4039.9333
$ws.Range("K122").Value = 12119.7999
$ws.Range("M122").Value = -9669.7999
$ws.Range("H133").Value = 89999.5
$ws.Range("J133").Value = 89999.5
$ws.Range("L133").Value = 89999.5
$ws.Range("N133").Value = -95059.5
$ws.Range("H136").Value = 3898.889
$ws.Range("I136").Value = 3755.5676
$ws.Range("J136").Value = 4561.75
$ws.Range("K136").Value = 11266.7028
$ws.Range("L136").Value = 13685.25
$ws.Range("M136").Value = -8716.702799999999
$ws.Range("N136").Value = -18785.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 27500
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H68").Value = 10000
$ws.Range("I68").Value = 10000
$ws.Range("K68").Value = 10000
$ws.Range("M68").Value = -9189
$ws.Range("H71").Value = 10000
$ws.Range("I71").Value = 10000
$ws.Range("K71").Value = 30000
$ws.Range("M71").Value = -25944
$ws.Range("H132").Value = 3580.1943
$ws.Range("I132").Value = 3159.68
$ws.Range("K132").Value = 9479.039999999999
$ws.Range("M132").Value = -6949.039999999999
$ws.Range("H135").Value = 83468.75
$ws.Range("J135").Value = 87187.5
$ws.Range("L135").Value = 87187.5
$ws.Range("N135").Value = -97327.5

Write-Host "Applied all Faerie_Profits updates across sheets."